# Add 2 more parts to Task 2 validation
# Appends two new worksheets ("results_3" and "results_4") at the end of
# the workbook, mirroring the layout/columns of the existing "results_*"
# sheets (Region, ElectricHeater_size_MMBtu_per_hr, Purchase_Price,
# Electricity_Price_per_MMBtu, Hourly_Cost, First_Year_Cost).

$wb = $excel.ActiveWorkbook

$headers = @(
    "Region",
    "ElectricHeater_size_MMBtu_per_hr",
    "Purchase_Price",
    "Electricity_Price_per_MMBtu",
    "Hourly_Cost",
    "First_Year_Cost"
)

$data3 = @(
    @("Midwest",   3.138, 50494.05, 20.35, 64.51, 391106.85),
    @("Northeast", 3.138, 50494.05, 24.47, 77.57, 460063.65),
    @("South",     3.138, 50494.05, 17.63, 55.88, 345540.45),
    @("West",      3.138, 50494.05, 24.09, 76.36, 453674.85)
)

$data4 = @(
    @("Midwest",   3.348, 53860.32, 20.35, 68.81, 417177.12),
    @("Northeast", 3.348, 53860.32, 24.47, 82.74, 490727.52),
    @("South",     3.348, 53860.32, 17.63, 59.61, 368601.12),
    @("West",      3.348, 53860.32, 24.09, 81.45, 483916.32)
)

# Add "results_3" right after the last existing sheet (resultsC_2).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "results_3"

# Add "results_4" right after "results_3".
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "results_4"

function Fill-ResultsSheet($ws, $headers, $rows) {
    for ($c = 1; $c -le $headers.Length; $c++) {
        $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
    }
    for ($r = 0; $r -lt $rows.Length; $r++) {
        $row = $rows[$r]
        for ($c = 0; $c -lt $row.Length; $c++) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
        }
    }
}

Fill-ResultsSheet $ws3 $headers $data3
Fill-ResultsSheet $ws4 $headers $data4
